{"js": "// Replace the date line and every \"AxB=C\" multiplication-table answer with\n// the new values from the commit. Each old value is unique in the document,\n// so a scoped search-and-replace (in document order) is unambiguous.\nconst replacements = [\n  [\"2025-08-04 Monday\", \"2025-08-05 Tuesday\"],\n  [\"58\u00d778=4524\", \"54\u00d764=3456\"],\n  [\"58\u00d781=4698\", \"30\u00d766=1980\"],\n  [\"81\u00d725=2025\", \"70\u00d789=6230\"],\n  [\"29\u00d770=2030\", \"21\u00d779=1659\"],\n  [\"85\u00d723=1955\", \"70\u00d765=4550\"],\n  [\"79\u00d799=7821\", \"49\u00d716=784\"],\n  [\"69\u00d763=4347\", \"89\u00d751=4539\"],\n  [\"90\u00d713=1170\", \"79\u00d750=3950\"],\n  [\"32\u00d731=992\", \"56\u00d750=2800\"],\n  [\"44\u00d770=3080\", \"39\u00d717=663\"],\n  [\"37\u00d733=1221\", \"66\u00d743=2838\"],\n  [\"93\u00d719=1767\", \"46\u00d792=4232\"],\n  [\"39\u00d720=780\", \"44\u00d730=1320\"],\n  [\"42\u00d746=1932\", \"25\u00d779=1975\"],\n  [\"90\u00d714=1260\", \"59\u00d752=3068\"],\n  [\"92\u00d751=4692\", \"17\u00d797=1649\"],\n  [\"97\u00d723=2231\", \"11\u00d776=836\"],\n  [\"55\u00d785=4675\", \"76\u00d767=5092\"],\n  [\"47\u00d760=2820\", \"48\u00d753=2544\"],\n  [\"82\u00d768=5576\", \"49\u00d788=4312\"],\n  [\"21\u00d795=1995\", \"56\u00d736=2016\"],\n  [\"74\u00d712=888\", \"81\u00d788=7128\"],\n  [\"18\u00d790=1620\", \"15\u00d752=780\"],\n  [\"62\u00d747=2914\", \"36\u00d725=900\"],\n  [\"28\u00d786=2408\", \"22\u00d715=330\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-08-04 Monday\", \"2025-08-05 Tuesday\"),\n  @(\"58\u00d778=4524\", \"54\u00d764=3456\"),\n  @(\"58\u00d781=4698\", \"30\u00d766=1980\"),\n  @(\"81\u00d725=2025\", \"70\u00d789=6230\"),\n  @(\"29\u00d770=2030\", \"21\u00d779=1659\"),\n  @(\"85\u00d723=1955\", \"70\u00d765=4550\"),\n  @(\"79\u00d799=7821\", \"49\u00d716=784\"),\n  @(\"69\u00d763=4347\", \"89\u00d751=4539\"),\n  @(\"90\u00d713=1170\", \"79\u00d750=3950\"),\n  @(\"32\u00d731=992\", \"56\u00d750=2800\"),\n  @(\"44\u00d770=3080\", \"39\u00d717=663\"),\n  @(\"37\u00d733=1221\", \"66\u00d743=2838\"),\n  @(\"93\u00d719=1767\", \"46\u00d792=4232\"),\n  @(\"39\u00d720=780\", \"44\u00d730=1320\"),\n  @(\"42\u00d746=1932\", \"25\u00d779=1975\"),\n  @(\"90\u00d714=1260\", \"59\u00d752=3068\"),\n  @(\"92\u00d751=4692\", \"17\u00d797=1649\"),\n  @(\"97\u00d723=2231\", \"11\u00d776=836\"),\n  @(\"55\u00d785=4675\", \"76\u00d767=5092\"),\n  @(\"47\u00d760=2820\", \"48\u00d753=2544\"),\n  @(\"82\u00d768=5576\", \"49\u00d788=4312\"),\n  @(\"21\u00d795=1995\", \"56\u00d736=2016\"),\n  @(\"74\u00d712=888\", \"81\u00d788=7128\"),\n  @(\"18\u00d790=1620\", \"15\u00d752=780\"),\n  @(\"62\u00d747=2914\", \"36\u00d725=900\"),\n  @(\"28\u00d786=2408\", \"22\u00d715=330\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $ok) {\n    throw \"Text not found: $old\"\n  }\n}\n"}
